$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 17-19 need the same style as existing data rows (bordered/bold A column)
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.749800258207028
$ws.Range("D10").Value = 0.2761562175175173
$ws.Range("E10").Value = 1.047988804147873
$ws.Range("F10").Value = 1.749800258207028
$ws.Range("G10").Value = 0.5946063655581701
$ws.Range("H10").Value = 1.16863442246823
$ws.Range("I10").Value = 1.13665739027051
$ws.Range("J10").Value = 0.2761562175175173
$ws.Range("K10").Value = 0.6620725108326952
$ws.Range("L10").Value = 1.205936384519862
$ws.Range("M10").Value = 0.9956405763615549

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.006228476626889
$ws.Range("D11").Value = 0.2779143930984984
$ws.Range("E11").Value = 1.111730827364047
$ws.Range("F11").Value = 1.006228476626889
$ws.Range("G11").Value = 1.164349911705717
$ws.Range("H11").Value = 0.7329516657898454
$ws.Range("I11").Value = 0.9417280633420867
$ws.Range("J11").Value = 0.2779143930984984
$ws.Range("K11").Value = 0.6948226102312728
$ws.Range("L11").Value = 0.8505255434290807
$ws.Range("M11").Value = 0.8724838896545138

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9898982363911495
$ws.Range("D12").Value = 0.2784406037713952
$ws.Range("E12").Value = 1.116686367963759
$ws.Range("F12").Value = 0.9898982363911495
$ws.Range("G12").Value = 1.166667449041278
$ws.Range("H12").Value = 0.7300729902747739
$ws.Range("I12").Value = 0.943046417429923
$ws.Range("J12").Value = 0.2784406037713952
$ws.Range("K12").Value = 0.6975634858675772
$ws.Range("L12").Value = 0.8437308611293634
$ws.Range("M12").Value = 0.8708020108120466

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.004171626794477
$ws.Range("D13").Value = 0.2778248932673483
$ws.Range("E13").Value = 1.114752660545838
$ws.Range("F13").Value = 1.004171626794477
$ws.Range("G13").Value = 1.164215448522576
$ws.Range("H13").Value = 0.7275662548484654
$ws.Range("I13").Value = 0.9425816035884211
$ws.Range("J13").Value = 0.2778248932673483
$ws.Range("K13").Value = 0.6962887769065933
$ws.Range("L13").Value = 0.8502302018505352
$ws.Range("M13").Value = 0.8718520812611876

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.0748960000000001
$ws.Range("D14").Value = 0.1761840000000006
$ws.Range("E14").Value = 1.196495999999997
$ws.Range("F14").Value = 0.0748960000000001
$ws.Range("G14").Value = 0.6403520000000004
$ws.Range("H14").Value = 5.223480000000011
$ws.Range("I14").Value = 0.8396599999999987
$ws.Range("J14").Value = 0.1761840000000006
$ws.Range("K14").Value = 0.686339999999999
$ws.Range("L14").Value = 0.3806179999999995
$ws.Range("M14").Value = 1.358511333333335

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.01
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1.017362500000002
$ws.Range("F15").Value = 0.01
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 8.854112499999994
$ws.Range("I15").Value = 0.7745874999999993
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0.5086812500000008
$ws.Range("L15").Value = 0.2593406250000004
$ws.Range("M15").Value = 1.776010416666666

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.4290100833279992
$ws.Range("D16").Value = 0.4300584056832036
$ws.Range("E16").Value = 1.047442264371203
$ws.Range("F16").Value = 0.4290100833279992
$ws.Range("G16").Value = 0.3965146263552006
$ws.Range("H16").Value = 5.379051837030389
$ws.Range("I16").Value = 0.8898836360192018
$ws.Range("J16").Value = 0.4300584056832036
$ws.Range("K16").Value = 0.7387503350272032
$ws.Range("L16").Value = 0.5838802091776012
$ws.Range("M16").Value = 1.428660142131199

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 1.002243512624621
$ws.Range("D17").Value = 0.9890388194373727
$ws.Range("E17").Value = 0.9828962381417687
$ws.Range("F17").Value = 1.002243512624621
$ws.Range("G17").Value = 0.9876195213177464
$ws.Range("H17").Value = 0.9967617647457098
$ws.Range("I17").Value = 0.9915358291192526
$ws.Range("J17").Value = 0.9890388194373727
$ws.Range("K17").Value = 0.9859675287895707
$ws.Range("L17").Value = 0.994105520707096
$ws.Range("M17").Value = 0.9916826142310787

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.7586173061555574
$ws.Range("D18").Value = 1.049502901720204
$ws.Range("E18").Value = 1.214131600674442
$ws.Range("F18").Value = 0.7586173061555574
$ws.Range("G18").Value = 1.047355794254368
$ws.Range("H18").Value = 1.04644598680168
$ws.Range("I18").Value = 0.9953858845190041
$ws.Range("J18").Value = 1.049502901720204
$ws.Range("K18").Value = 1.131817251197323
$ws.Range("L18").Value = 0.9452172786764403
$ws.Range("M18").Value = 1.018573245687543

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9925697688517625
$ws.Range("D19").Value = 1.403486764409819
$ws.Range("E19").Value = 0.8934585392115453
$ws.Range("F19").Value = 0.9925697688517625
$ws.Range("G19").Value = 1.119167319757488
$ws.Range("H19").Value = 0.7320442331324901
$ws.Range("I19").Value = 0.9211748398454243
$ws.Range("J19").Value = 1.403486764409819
$ws.Range("K19").Value = 1.148472651810682
$ws.Range("L19").Value = 1.070521210331222
$ws.Range("M19").Value = 1.010316910868088
